# Add two new element/type rows ("solar_th" and "pvt") to the table,
# appended right after the existing "CHP1"/"CHP" row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "solar_th1"
$ws.Range("B7").Value = "solar_th"

$ws.Range("A8").Value = "pvt1"
$ws.Range("B8").Value = "pvt"
